# Auto-generated Excel COM-interop script applying the "Update latest output (run 23)" diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" (sheet1) ----
$schedule = $wb.Worksheets.Item("Schedule")

# Row 2
$schedule.Cells.Item(2,1).Value = 46037.04166666666
$schedule.Cells.Item(2,3).Value = 4
$schedule.Cells.Item(2,4).Value = 15.12
$schedule.Cells.Item(2,5).Value = 587.08463775
$schedule.Cells.Item(2,6).Value = 38.82834905753968

# Row 4
$schedule.Cells.Item(4,1).Value = 46037.875
$schedule.Cells.Item(4,2).Value = 46038.04166666666
$schedule.Cells.Item(4,3).Value = 4
$schedule.Cells.Item(4,4).Value = 15.12
$schedule.Cells.Item(4,5).Value = 590.1241807499999
$schedule.Cells.Item(4,6).Value = 39.02937703373016

# Row 5
$schedule.Cells.Item(5,1).Value = 46038.375
$schedule.Cells.Item(5,3).Value = 7
$schedule.Cells.Item(5,4).Value = 26.46
$schedule.Cells.Item(5,5).Value = 468.1814669999999
$schedule.Cells.Item(5,6).Value = 17.69393299319728

# New row 6 (extends the schedule with an additional pump-on interval)
$schedule.Cells.Item(6,1).Value = 46038.83333333334
$schedule.Cells.Item(6,2).Value = 46039
$schedule.Cells.Item(6,3).Value = 4
$schedule.Cells.Item(6,4).Value = 15.12
$schedule.Cells.Item(6,5).Value = 400.345803
$schedule.Cells.Item(6,6).Value = 26.47789702380953

# Copy style/number-format from row 5 so the new date cells (A6:B6) render as datetimes
$schedule.Range("A6:B6").Style = $schedule.Range("A5:B5").Style
$schedule.Range("A6:B6").NumberFormat = $schedule.Range("A5:B5").NumberFormat

# ---- Sheet "Detailed" (sheet2) ----
$detailed = $wb.Worksheets.Item("Detailed")

# Row 3
$detailed.Cells.Item(3,5).Value = "OFF"

# Row 43
$detailed.Cells.Item(43,2).Value = 120.01

# Row 44
$detailed.Cells.Item(44,2).Value = 85.95
$detailed.Cells.Item(44,5).Value = "ON"

# Row 45
$detailed.Cells.Item(45,2).Value = 80.02
$detailed.Cells.Item(45,3).Value = "historical"

# Row 46
$detailed.Cells.Item(46,2).Value = 85.95
$detailed.Cells.Item(46,3).Value = "historical"

# Row 47
$detailed.Cells.Item(47,2).Value = 71.40000000000001

# Row 48
$detailed.Cells.Item(48,2).Value = 65

# Row 49
$detailed.Cells.Item(49,2).Value = 73.93557

# Row 50
$detailed.Cells.Item(50,2).Value = 65

# Row 51
$detailed.Cells.Item(51,2).Value = 78

# Row 52
$detailed.Cells.Item(52,2).Value = 58.69504
$detailed.Cells.Item(52,5).Value = "OFF"

# Row 53
$detailed.Cells.Item(53,2).Value = 58.45285
$detailed.Cells.Item(53,5).Value = "OFF"

# Row 54
$detailed.Cells.Item(54,2).Value = 59.60656
$detailed.Cells.Item(54,5).Value = "OFF"

# Row 55
$detailed.Cells.Item(55,2).Value = 64.89
$detailed.Cells.Item(55,5).Value = "OFF"

# Row 56
$detailed.Cells.Item(56,2).Value = 64.89
$detailed.Cells.Item(56,5).Value = "OFF"

# Row 57
$detailed.Cells.Item(57,2).Value = 58.69131
$detailed.Cells.Item(57,5).Value = "OFF"

# Row 58
$detailed.Cells.Item(58,2).Value = 64.89

# Row 59
$detailed.Cells.Item(59,2).Value = 80.02

# Row 60
$detailed.Cells.Item(60,2).Value = 80.02

# Row 61
$detailed.Cells.Item(61,2).Value = 85.95

# Row 62
$detailed.Cells.Item(62,2).Value = 79.95

# Row 64
$detailed.Cells.Item(64,2).Value = 36.07

# Row 65
$detailed.Cells.Item(65,2).Value = 48.49273

# Row 66
$detailed.Cells.Item(66,2).Value = 57.06003
$detailed.Cells.Item(66,5).Value = "OFF"

# Row 67
$detailed.Cells.Item(67,2).Value = 55.27342
$detailed.Cells.Item(67,5).Value = "OFF"

# Row 72
$detailed.Cells.Item(72,2).Value = 36.05949

# Row 73
$detailed.Cells.Item(73,2).Value = 36.05989

# Row 79
$detailed.Cells.Item(79,2).Value = 33.07075

# Row 80
$detailed.Cells.Item(80,2).Value = 27.77113

# Row 81
$detailed.Cells.Item(81,2).Value = 22.68446

# Row 83
$detailed.Cells.Item(83,2).Value = 0.81985

# Row 84
$detailed.Cells.Item(84,2).Value = -9.431789999999999

# Row 85
$detailed.Cells.Item(85,2).Value = -8.136850000000001

# Row 86
$detailed.Cells.Item(86,2).Value = -6.78314

# Row 87
$detailed.Cells.Item(87,2).Value = -3.00909

# Row 88
$detailed.Cells.Item(88,2).Value = -3.07152

# Row 89
$detailed.Cells.Item(89,2).Value = 29.85322

# Row 90
$detailed.Cells.Item(90,2).Value = 32.40461
$detailed.Cells.Item(90,5).Value = "ON"

# Row 91
$detailed.Cells.Item(91,5).Value = "ON"

# Row 92
$detailed.Cells.Item(92,2).Value = 29.85322
$detailed.Cells.Item(92,5).Value = "ON"

# Row 93
$detailed.Cells.Item(93,5).Value = "ON"

# Row 94
$detailed.Cells.Item(94,5).Value = "ON"

# Row 95
$detailed.Cells.Item(95,5).Value = "ON"

# Row 96
$detailed.Cells.Item(96,2).Value = 58.94929
$detailed.Cells.Item(96,5).Value = "ON"

# Row 97
$detailed.Cells.Item(97,5).Value = "ON"

